# Add a new "2022-Q4" worksheet (fund holding detail for 002153 / 石基信息)
# right after the "总计" (summary) sheet, and record its summary stats on
# the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying the "2022-Q3" sheet (same
#    layout/styles as every other quarterly detail sheet) and placing the
#    copy right after "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $wb.Worksheets.Item("总计"))
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Wipe the copied (2022-Q3) values - keep formatting/styles intact.
$newSheet.Range("A2:H13").ClearContents()

# The fund-code / amount / position columns (B:G) are stored as TEXT in
# this workbook (so leading zeros in fund codes like "010115" survive).
# Force text formatting before writing so Excel doesn't coerce these
# numeric-looking strings into numbers.
$newSheet.Range("B2:G29").NumberFormat = "@"

# Column A keeps the bold/centered "index" style used by row 2-13 (s=2 in
# the template); extend that same formatting down through row 29.
$newSheet.Range("A2").Copy($newSheet.Range("A14:A29"))

$data2022Q4 = @(
  @(0, '010115', '易方达远见成长混合A', '28.67', '91.44', '5.25', '1.5052', 3),
  @(1, '000118', '广发聚鑫债券A', '123.74', '20.30', '0.92', '1.1384', 9),
  @(2, '011412', '易方达远见成长混合C', '18.27', '91.44', '5.25', '0.9592', 3),
  @(3, '001437', '易方达瑞享灵活配置混合I', '10.56', '92.06', '6.94', '0.7329', 4),
  @(4, '001438', '易方达瑞享灵活配置混合E', '10.56', '92.06', '6.94', '0.7329', 4),
  @(5, '410003', '华富成长趋势混合', '10.74', '87.35', '5.14', '0.5520', 7),
  @(6, '009121', '广发招享混合A', '52.83', '27.71', '0.98', '0.5177', 10),
  @(7, '011891', '易方达先锋成长混合A', '7.13', '92.85', '7.10', '0.5062', 2),
  @(8, '410007', '华富价值增长混合', '8.09', '79.81', '4.78', '0.3867', 8),
  @(9, '011892', '易方达先锋成长混合C', '4.41', '92.85', '7.10', '0.3131', 2),
  @(10, '501062', '南方瑞合三年定期开放混合（LOF）', '7.50', '89.25', '4.10', '0.3075', 7),
  @(11, '012408', '广发恒昌一年持有期混合A', '24.13', '27.85', '1.04', '0.2510', 7),
  @(12, '013880', '广发招享混合C', '23.13', '27.71', '0.98', '0.2267', 10),
  @(13, '006864', '国联安核心资产策略混合', '4.67', '91.48', '3.59', '0.1677', 9),
  @(14, '009398', '华富成长企业精选股票', '3.12', '94.55', '5.19', '0.1619', 8),
  @(15, '000119', '广发聚鑫债券C', '15.35', '20.30', '0.92', '0.1412', 9),
  @(16, '012586', '南方港股创新视野一年持有混合A', '2.20', '85.91', '3.73', '0.0821', 9),
  @(17, '014706', '华富匠心明选一年持有期混合A', '1.93', '89.78', '3.91', '0.0755', 10),
  @(18, '010925', '兴银科技增长1个月滚动持有期混合A', '0.96', '81.72', '7.73', '0.0742', 1),
  @(19, '519644', '银河智联主题灵活配置混合', '1.20', '89.42', '4.89', '0.0587', 7),
  @(20, '014707', '华富匠心明选一年持有期混合C', '1.48', '89.78', '3.91', '0.0579', 10),
  @(21, '015412', '西部利得数字产业混合A', '0.94', '92.30', '5.12', '0.0481', 6),
  @(22, '003152', '华富天鑫灵活配置混合A', '0.93', '88.63', '4.56', '0.0424', 7),
  @(23, '015413', '西部利得数字产业混合C', '0.43', '92.30', '5.12', '0.0220', 6),
  @(24, '012409', '广发恒昌一年持有期混合C', '1.99', '27.85', '1.04', '0.0207', 7),
  @(25, '010926', '兴银科技增长1个月滚动持有期混合C', '0.11', '81.72', '7.73', '0.0085', 1),
  @(26, '003153', '华富天鑫灵活配置混合C', '0.16', '88.63', '4.56', '0.0073', 7),
  @(27, '012587', '南方港股创新视野一年持有混合C', '0.17', '85.91', '3.73', '0.0063', 9)
)

$r = 2
foreach ($row in $data2022Q4) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Drop the transient "text number format" now that the values are safely
# stored as text - restores the plain/default style the other detail
# sheets use for their data rows (B:G carry no explicit style).
$newSheet.Range("B2:G29").Style = "Normal"

# ---------------------------------------------------------------------
# 2. Record the 2022-Q4 summary row on "总计": 28 holdings, 9.1 亿元.
#    Shift the existing rows (2022-Q3 .. 2020-Q4) down by one row first.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
for ($r = 9; $r -ge 2; $r--) {
    $src = $summary.Range("A" + $r + ":D" + $r)
    $dst = $summary.Range("A" + ($r + 1) + ":D" + ($r + 1))
    $src.Copy($dst)
}

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 28
$summary.Range("D2").Value = 9.1

# ---------------------------------------------------------------------
# 3. Keep "2020-Q4" as the active/selected tab (matches the workbook's
#    original selection).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
